# Apply the "Add files via upload" header-restructuring edit to
# ImportStaff.xlsx's "Staff Vitals" and "Staff Style" sheets.
#
# NOTE: this runtime's EntireColumn.Delete() does not reliably honor a
# multi-area Range (e.g. "D1,E1,I1,K1") - it only removes one column
# worth of width. So every column deletion below is issued individually,
# rightmost column first, so earlier (left-of-it) column letters stay
# valid for the subsequent deletes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Staff Vitals": add Face ID / Height / Position in front, drop
# the old FACEID and UNIQUE_PHOTO_ID columns, and append a renamed
# "*Unique Photo ID" column at the end.
# ---------------------------------------------------------------------
$wsVitals = $wb.Worksheets.Item("Staff Vitals")

# Remove "Staff Vitals - UNIQUE_PHOTO_ID" (U) and "Staff Vitals - FACEID" (F).
# Delete the rightmost one first so column "F" is still correct afterwards.
$wsVitals.Range("U1").EntireColumn.Delete()
$wsVitals.Range("F1").EntireColumn.Delete()

# Insert 3 blank columns at the front for the new Face ID / Height / Position fields.
$wsVitals.Range("A1:C1").EntireColumn.Insert()

# Give the new header cells the same formatting (bold, border, alignment) as
# the existing header row before filling in their text.
$wsVitals.Range("D1").Copy()
$wsVitals.Range("A1:C1").PasteSpecial(-4122)

$wsVitals.Range("A1").Value() = "Face ID"
$wsVitals.Range("B1").Value() = "Height"
$wsVitals.Range("C1").Value() = "Position"

# Append the renamed unique-photo-id column at the new end of the row (Z1).
$wsVitals.Range("Y1").Copy()
$wsVitals.Range("Z1").PasteSpecial(-4122)
$wsVitals.Range("Z1").Value() = "*Unique Photo ID"

# ---------------------------------------------------------------------
# Sheet "Staff Style": rename the first header, drop four old
# proficiency columns, and insert five friendlier-named proficiency
# headers after the first column.
# ---------------------------------------------------------------------
$wsStyle = $wb.Worksheets.Item("Staff Style")

# Remove "Staff Style - DEFENSE_PROFICIENCY" (D), "...GRIT_&_GRIND_PROFICIENCY" (E),
# "...PACE_&_SPACE_PROFICIENCY" (I) and "...POST_CENTRIC_PROFICIENCY" (K).
# Delete rightmost first so earlier letters remain valid.
$wsStyle.Range("K1").EntireColumn.Delete()
$wsStyle.Range("I1").EntireColumn.Delete()
$wsStyle.Range("E1").EntireColumn.Delete()
$wsStyle.Range("D1").EntireColumn.Delete()

# Rename the first header in place.
$wsStyle.Range("A1").Value() = "Seven Seconds Proficiency"

# Insert 5 new columns right after A for the new proficiency headers.
$wsStyle.Range("B1:F1").EntireColumn.Insert()

# Copy formatting from an existing header cell into the new ones.
$wsStyle.Range("G1").Copy()
$wsStyle.Range("B1:F1").PasteSpecial(-4122)

$wsStyle.Range("B1").Value() = "Defense Proficiency"
$wsStyle.Range("C1").Value() = "Grit & Grind Proficiency"
$wsStyle.Range("D1").Value() = "Pace &Space Proficiency"
$wsStyle.Range("E1").Value() = "Perimeter Centric Proficiency"
$wsStyle.Range("F1").Value() = "Post Centric Proficiency"
